$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, splitting the old "Condition"/"Replicate"
# column pair into two marker columns: "Comparison1" (old C) and "Comparison2" (new D).
# Everything that used to live in columns D:N shifts right to E:O.
$ws.Columns.Item(4).Insert()

# The old column C held the numeric "Replicate" value (1 or 2) per row. Clear it out;
# it will be replaced below by an "x" marker in either the new C or D column.
$ws.Range("C2:C5").ClearContents()

# Mark which comparison each row belongs to with an "x" (set data cells before the
# header text so the shared-string table ends up in the same append order as the
# authored workbook: ... "<copy path here>", "x", "Comparison1", "Comparison2").
$ws.Range("C2").Value = "x"
$ws.Range("D3").Value = "x"
$ws.Range("C4").Value = "x"
$ws.Range("D5").Value = "x"

# New headers for the split comparison columns.
$ws.Range("C1").Value = "Comparison1"
$ws.Range("D1").Value = "Comparison2"

# Match the widened header columns (approx. character width of the new "Comparison1"/
# "Comparison2" headers).
$ws.Range("C1:D1").ColumnWidth = 11.7

# Move the active selection to D4, matching the post-edit cursor position.
$ws.Range("D4").Select()
